$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.626.35"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "3.327.51"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "587.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.131"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.412"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "3.915.56"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "68.763.43"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000169"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "3.320.49"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "447.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.521"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "3.484.80"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.190"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +5.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.794"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("D43").Value = "2.693.71"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0682"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "329.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0279"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.82%  "
